{"js": "// The document has a centered title paragraph with a date/weekday line,\n// followed by a single table of two-digit \u00f7 one-digit division problems.\n// Every row has 5 cells, but only every 4th row (0, 4, 8, 12, 16, ...)\n// actually holds text - the rows in between are blank spacer rows with\n// empty cells. The edit replaces the title text and every non-empty\n// problem cell's text, in document (row-major) order, with new values.\n// Some new values coincide with *old* values that live elsewhere in the\n// table, so cells must be updated by position in a single pass - not\n// via find/replace (which would risk re-matching an already-updated\n// cell or a not-yet-updated one).\n\nconst newTitle = \"2025-02-11 Tuesday\";\n\n// New cell text, in row-major reading order (only the filled rows,\n// 5 cells each = 25 values total).\nconst newValues = [\n  \"65\u00f79=7, 2\", \"17\u00f72=8, 1\", \"80\u00f78=10, 0\", \"87\u00f77=12, 3\", \"23\u00f73=7, 2\",\n  \"45\u00f72=22, 1\", \"20\u00f75=4, 0\", \"77\u00f75=15, 2\", \"73\u00f76=12, 1\", \"25\u00f72=12, 1\",\n  \"69\u00f73=23, 0\", \"85\u00f78=10, 5\", \"76\u00f76=12, 4\", \"65\u00f72=32, 1\", \"99\u00f78=12, 3\",\n  \"45\u00f78=5, 5\", \"43\u00f78=5, 3\", \"27\u00f79=3, 0\", \"99\u00f75=19, 4\", \"12\u00f79=1, 3\",\n  \"78\u00f74=19, 2\", \"29\u00f73=9, 2\", \"81\u00f79=9, 0\", \"34\u00f74=8, 2\", \"69\u00f72=34, 1\"\n];\n\nconst body = context.document.body;\n\n// --- Update the title paragraph (first paragraph in the body). ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(newTitle, Word.InsertLocation.replace);\n\n// --- Update every non-blank problem cell in the (single) table. ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nlet valueIndex = 0;\nfor (let r = 0; r < table.rowCount && valueIndex < newValues.length; r++) {\n  const rowValues = table.values[r];\n  for (let c = 0; c < rowValues.length && valueIndex < newValues.length; c++) {\n    if (rowValues[c] === \"\") continue;\n    table.getCell(r, c).value = newValues[valueIndex];\n    valueIndex++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document has a centered title paragraph with a date/weekday line,\n# followed by a single table of two-digit / one-digit division problems.\n# Every row of the table has 5 cells, but only every 4th row (1, 5, 9,\n# 13, 17 in 1-based Word numbering) actually holds text - the rows in\n# between are blank spacer rows with empty cells. The edit replaces the\n# title text and every non-empty problem cell's text, in document\n# (row-major) order, with new values. Some new values coincide with\n# *old* values that live elsewhere in the table, so cells must be\n# updated by position in a single pass - not via Find/Replace (which\n# would risk re-matching an already-updated cell or a not-yet-updated\n# one).\n\n$d = $word.ActiveDocument\n\n$newTitle = \"2025-02-11 Tuesday\"\n\n# New cell text, in row-major reading order (only the filled rows,\n# 5 cells each = 25 values total).\n$newValues = @(\n  \"65\u00f79=7, 2\", \"17\u00f72=8, 1\", \"80\u00f78=10, 0\", \"87\u00f77=12, 3\", \"23\u00f73=7, 2\",\n  \"45\u00f72=22, 1\", \"20\u00f75=4, 0\", \"77\u00f75=15, 2\", \"73\u00f76=12, 1\", \"25\u00f72=12, 1\",\n  \"69\u00f73=23, 0\", \"85\u00f78=10, 5\", \"76\u00f76=12, 4\", \"65\u00f72=32, 1\", \"99\u00f78=12, 3\",\n  \"45\u00f78=5, 5\", \"43\u00f78=5, 3\", \"27\u00f79=3, 0\", \"99\u00f75=19, 4\", \"12\u00f79=1, 3\",\n  \"78\u00f74=19, 2\", \"29\u00f73=9, 2\", \"81\u00f79=9, 0\", \"34\u00f74=8, 2\", \"69\u00f72=34, 1\"\n)\n\n# --- Update the title paragraph (first paragraph in the body). ---\n$d.Paragraphs.Item(1).Range.Text = $newTitle\n\n# --- Update every non-blank problem cell in the (single) table. ---\n$t = $d.Tables.Item(1)\n\n$valueIndex = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        if ($valueIndex -ge $newValues.Count) { break }\n\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text\n        # Strip the trailing cell-mark character(s) Word appends to\n        # Range.Text (\"\\r\\a\") before checking for blankness.\n        $plain = $cellText.TrimEnd([char]13, [char]7)\n\n        if ($plain -eq \"\") { continue }\n\n        $cell.Range.Text = $newValues[$valueIndex]\n        $valueIndex++\n    }\n}\n"}
